# CR_itemshape_2: tidy up the WG_Number_excel_table sheet.
#
# - Row 31 was a stray/duplicate leftover row (only F/G populated, holding
#   the cancelled "N9358" entry and a "CANCELLED ... ARM EXPRESS" note) that
#   doesn't belong in the clean N-number sequence (N9357, N9358(removed),
#   N9359, ...) - delete it so the table goes straight from N9357 to N9359
#   and everything below shifts up one row.
# - Column G held ad-hoc "Changed from edX to edY" / cancellation comments
#   that are no longer needed - delete the whole column.
# - Reset the view: scroll back to the top-left corner (A1) instead of the
#   old scrolled position, and leave the selection on E19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the active window back to the top-left (A1) - this is what drops
# the old topLeftCell="A10" setting.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1

# Remove the stray row 31 (N9358 / CANCELLED note); rows below shift up.
$ws.Rows("31").Delete()

# Remove column G (the "Changed from ed.. to ed.." / CANCELLED comments column).
$ws.Columns("G").Delete()

# Match the recorded selection in the saved view.
$ws.Range("E19").Select()
